$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the two template placeholders in the control-part unit-cost block.
# K3 previously held ${data.UNIT_TM_AMT}; it now reflects the merged
# "material finish" TM amount field.
$ws.Cells.Item(3, 11).Value = '${data.UNIT_MATERIAL_FINISH_TM_AMT}'

# K4 previously held ${data.UNIT_HEAT_AMT}; it now reflects the merged
# "material finish" HEAT amount field.
$ws.Cells.Item(4, 11).Value = '${data.UNIT_MATERIAL_FINISH_HEAT_AMT}'
